$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.02672833333333334
$ws.Range("H2").Value = 0.08018500000000001
$ws.Range("I2").Value = 0.3128877685602129
$ws.Range("J2").Value = 0.3128877685602129
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 13.69681033333333
$ws.Range("N2").Value = 41.090431
$ws.Range("O2").Value = 0.1107101339353595
$ws.Range("P2").Value = 0.1107101339353595
$ws.Range("Q2").Value = 0.3660929121927778
$ws.Range("R2").Value = 3.294836209735
$ws.Range("S2").Value = 0.03463984676403695
$ws.Range("T2").Value = 0.03463984676403695

# Row 3
$ws.Range("G3").Value = 0.02672833333333334
$ws.Range("H3").Value = 0.08018500000000001
$ws.Range("I3").Value = 0.3128877685602129
$ws.Range("J3").Value = 0.3128877685602129
$ws.Range("O3").Value = 0.8037307792188669
$ws.Range("P3").Value = 0.803730779218867
$ws.Range("Q3").Value = 2.657752557277222
$ws.Range("R3").Value = 23.919773015495
$ws.Range("S3").Value = 0.2514775300329524
$ws.Range("T3").Value = 0.2514775300329525

# Row 4
$ws.Range("G4").Value = 0.02672833333333334
$ws.Range("H4").Value = 0.08018500000000001
$ws.Range("I4").Value = 0.3128877685602129
$ws.Range("J4").Value = 0.3128877685602129
$ws.Range("O4").Value = 0.08555908684577355
$ws.Range("P4").Value = 0.08555908684577354
$ws.Range("Q4").Value = 0.2829241927050001
$ws.Range("R4").Value = 2.546317734345001
$ws.Range("S4").Value = 0.02677039176322355
$ws.Range("T4").Value = 0.02677039176322355

# Row 5
$ws.Range("G5").Value = 0.05869633333333333
$ws.Range("I5").Value = 0.6871122314397871
$ws.Range("J5").Value = 0.6871122314397871
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 13.69681033333333
$ws.Range("N5").Value = 41.090431
$ws.Range("O5").Value = 0.1107101339353595
$ws.Range("P5").Value = 0.1107101339353595
$ws.Range("Q5").Value = 0.8039525449287777
$ws.Range("R5").Value = 7.235572904358999
$ws.Range("S5").Value = 0.07607028717132257
$ws.Range("T5").Value = 0.07607028717132257

# Row 6
$ws.Range("G6").Value = 0.05869633333333333
$ws.Range("I6").Value = 0.6871122314397871
$ws.Range("J6").Value = 0.6871122314397871
$ws.Range("O6").Value = 0.8037307792188669
$ws.Range("P6").Value = 0.803730779218867
$ws.Range("Q6").Value = 5.836515433789221
$ws.Range("S6").Value = 0.5522532491859145
$ws.Range("T6").Value = 0.5522532491859146

# Row 7
$ws.Range("G7").Value = 0.05869633333333333
$ws.Range("I7").Value = 0.6871122314397871
$ws.Range("J7").Value = 0.6871122314397871
$ws.Range("O7").Value = 0.08555908684577355
$ws.Range("P7").Value = 0.08555908684577354
$ws.Range("R7").Value = 5.591800754793001
$ws.Range("S7").Value = 0.05878869508255
$ws.Range("T7").Value = 0.05878869508254999
